$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 86, pushing existing rows 86:91 down to 89:94.
$ws.Rows.Item(86).Resize(3).Insert()

# Data for the 3 new rows (a new weekly report date: 2021-11-09 / serial 44509)
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Categoría ID,
#          G Categoría, H Variedad, I Calidad, J Volumen, K Precio minimo,
#          L Precio maximo, M Precio promedio ponderado, N Unidad de comercializacion,
#          O Origen, P Precio $/Kg, Q Kg o Unidades, R Clasificacion

$newRows = @(
    @{ Row = 86; I = "Banquete"; J = 1950; K = 1200; L = 1300; M = 1218; N = '$/paquete'; O = "Provincia de Linares"; P = 1218 },
    @{ Row = 87; I = "Primera";  J = 1850; K = 1000; L = 1100; M = 1024; N = '$/paquete'; O = "Provincia de Linares"; P = 1024 },
    @{ Row = 88; I = "Segunda";  J = 1250; K = 800;  L = 900;  M = 820;  N = '$/paquete'; O = "Provincia de Linares"; P = 820  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44509
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 300000000
    $ws.Cells.Item($row, 7).Value = "Espárragos"
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
